# Applies the cryptos list refresh described by the commit message
# "Updated cryptos list on Sat Jun 15 09:35:42 UTC 2024 with GitHub Actions".
#
# Every row's Price (column D) and Volume(1h) (column E) text is refreshed
# with newer scraped figures. Rows 10 and 11 also swap rankings (Toncoin
# moves up to rank 10, Dogecoin drops to rank 11), so their Coin name, Link
# and Price/Volume values move together.
#
# All Price/Volume cells in this sheet are stored as plain text (not
# numbers), even when a price happens to look like a clean decimal (e.g.
# "607.97"). Assigning such a string straight to .Value would make Excel
# auto-convert it into a genuine number, so those particular cells are
# written with a leading apostrophe (forces text entry, same as typing it
# in the Excel UI) and then reset to the default "Normal" style so no
# stray text-number-format style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.196.89"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "3.525.82"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'607.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'143.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D7").Value = "3.522.94"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D10").Value = "'8.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.137"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.52%  "
$ws.Range("E12").Value = "  -2.84%  "
$ws.Range("D13").Value = "4.126.11"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  -5.19%  "
$ws.Range("D15").Value = "'30.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.51%  "
$ws.Range("D16").Value = "3.530.35"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").Value = "66.314.51"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "'10.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").Value = "'6.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("D21").Value = "'14.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.57%  "
$ws.Range("D22").Value = "'425.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").Value = "'0.601"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.52%  "
$ws.Range("D24").Value = "'78.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "3.670.18"
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'0.0000119"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "'9.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.86%  "
$ws.Range("D29").Value = "'8.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.83%  "
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("E33").Value = "  -6.51%  "
$ws.Range("D34").Value = "'25.27"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").Value = "3.519.94"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("E37").Value = "  -3.25%  "
$ws.Range("D38").Value = "'7.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.71%  "
$ws.Range("E39").Value = "  -5.87%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D41").Value = "'172.59"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.28%  "
$ws.Range("E42").Value = "  -4.44%  "
$ws.Range("D43").Value = "'5.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.02%  "
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("E45").Value = "  -9.43%  "
$ws.Range("D46").Value = "'45.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "'25.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.73%  "
$ws.Range("E48").Value = "  -6.42%  "
$ws.Range("E49").Value = "  -2.29%  "
$ws.Range("D50").Value = "'7.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.945"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.64%  "

